$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.203.85"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "1.855.48"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'241.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("D6").Value = "'0.6982"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.07765"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.3070"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "

# Row 10
$ws.Range("D10").Value = "'23.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
$ws.Range("D11").Value = "'0.07813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.07%  "

# Row 12
$ws.Range("D12").Value = "1.852.12"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13
$ws.Range("E13").Value = "  -1.24%  "

# Row 14
$ws.Range("D14").Value = "'92.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.48%  "

# Row 15
$ws.Range("D15").Value = "'0.6864"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16
$ws.Range("D16").Value = "'6.517"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.72%  "

# Row 17
$ws.Range("D17").Value = "'0.000008472"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.45%  "

# Row 18
$ws.Range("D18").Value = "29.208.34"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19
$ws.Range("D19").Value = "'248.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "

# Row 20
$ws.Range("D20").Value = "2.107.00"
$ws.Range("E20").Value = "  -0.16%  "

# Row 21
$ws.Range("D21").Value = "'12.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.88%  "

# Row 22
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").Value = "'7.521"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.34%  "

# Row 24
$ws.Range("D24").Value = "'0.9997"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

# Row 25
$ws.Range("D25").Value = "'0.1496"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.24%  "

# Row 26
$ws.Range("D26").Value = "'161.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "

# Row 27
$ws.Range("D27").Value = "'8.859"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "

# Row 28
$ws.Range("D28").Value = "'18.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.55%  "

# Row 29
$ws.Range("D29").Value = "'1.553"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.92%  "

# Row 30
$ws.Range("D30").Value = "'4.244"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.08%  "

# Row 31
$ws.Range("D31").Value = "'4.207"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.06%  "

# Row 32
$ws.Range("E32").Value = "  -0.92%  "

# Row 33
$ws.Range("D33").Value = "'0.05219"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.48%  "

# Row 34
$ws.Range("D34").Value = "'0.7592"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "

# Row 35
$ws.Range("D35").Value = "'1.170"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.52%  "

# Row 36
$ws.Range("D36").Value = "'1.841"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.86%  "

# Row 37
$ws.Range("D37").Value = "'2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("D38").Value = "'0.01862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "

# Row 39
$ws.Range("D39").Value = "1.224.77"
$ws.Range("E39").Value = "  -1.42%  "

# Row 40
$ws.Range("D40").Value = "'2.727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "

# Row 41
$ws.Range("D41").Value = "'0.8989"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "

# Row 42
$ws.Range("D42").Value = "'109.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "

# Row 43
$ws.Range("E43").Value = "  -0.10%  "

# Row 44
$ws.Range("D44").Value = "'5.515"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.14%  "

# Row 45
$ws.Range("D45").Value = "2.005.14"
$ws.Range("E45").Value = "  -0.20%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'65.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.60%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000123"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.97%  "

# Row 48
$ws.Range("D48").Value = "'0.5180"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("D49").Value = "'9.520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.23%  "

# Row 50
$ws.Range("D50").Value = "'1.752"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "

# Row 51
$ws.Range("D51").Value = "'7.042"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.70%  "
